$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.21"
$ws.Range("E2").Value = "'0.88%"
$ws.Range("D3").Value = "'41.43"
$ws.Range("E3").Value = "'4.81%"
$ws.Range("D4").Value = "'5.645"
$ws.Range("E4").Value = "'-1.04%"
$ws.Range("D5").Value = "'0.08293"
$ws.Range("E5").Value = "'3.12%"
$ws.Range("D6").Value = "'2.047"
$ws.Range("E6").Value = "'1.15%"
$ws.Range("D7").Value = "'8.768"
$ws.Range("E7").Value = "'1.54%"
$ws.Range("E8").Value = "'0.83%"
$ws.Range("D9").Value = "'2.974"
$ws.Range("E9").Value = "'1.06%"
$ws.Range("D10").Value = "'0.9251"
$ws.Range("E10").Value = "'0.16%"
$ws.Range("D11").Value = "'0.1271"
$ws.Range("E11").Value = "'1.05%"
$ws.Range("D12").Value = "'0.1960"
$ws.Range("E12").Value = "'-0.49%"
$ws.Range("D13").Value = "'0.09355"
$ws.Range("E13").Value = "'1.82%"
$ws.Range("D14").Value = "'0.03932"
$ws.Range("E14").Value = "'10.24%"
$ws.Range("E15").Value = "'0.94%"
$ws.Range("D16").Value = "'0.001309"
$ws.Range("E16").Value = "'0.95%"
$ws.Range("D17").Value = "'0.006147"
$ws.Range("E17").Value = "'-3.07%"
$ws.Range("D19").Value = "'3.445"
$ws.Range("D21").Value = "'8.349"
$ws.Range("E21").Value = "'-4.64%"
$ws.Range("E22").Value = "'1.75%"
$ws.Range("E23").Value = "'-1.36%"
$ws.Range("D24").Value = "'0.04412"
$ws.Range("E24").Value = "'0.42%"
$ws.Range("D25").Value = "'0.001258"
$ws.Range("E25").Value = "'-0.25%"
$ws.Range("E26").Value = "'-6.41%"
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("E27").Value = "'1.06%"
$ws.Range("D39").Value = "'0.02782"
$ws.Range("E39").Value = "'11.38%"
$ws.Range("D40").Value = "'0.05512"
$ws.Range("E40").Value = "'3.45%"
$ws.Range("D41").Value = "'0.007938"
$ws.Range("E41").Value = "'6.50%"
$ws.Range("D42").Value = "'0.1423"
$ws.Range("E42").Value = "'1.20%"
$ws.Range("D43").Value = "'0.008938"
$ws.Range("E43").Value = "'-9.75%"
$ws.Range("D44").Value = "'0.002141"
$ws.Range("E44").Value = "'1.37%"
$ws.Range("E45").Value = "'8.40%"
$ws.Range("D46").Value = "'0.00007016"
$ws.Range("E46").Value = "'5.15%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.11%"
$ws.Range("D48").Value = "'0.003197"
$ws.Range("E48").Value = "'5.27%"
$ws.Range("E49").Value = "'0.09%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.11%"
